# Improved Tristan text, fixed chests in gatekeeper's house
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Todo sheet: remove the "Add gatekeeper chests (or fix them)" item
#    (the chests have now been fixed), leaving the row empty just like
#    the existing gap at row 7.
# ---------------------------------------------------------------------
$todo = $wb.Worksheets.Item("Todo")
$todo.Range("A3").ClearContents()

# ---------------------------------------------------------------------
# 2. Maps sheet: add a new map entry for the Gatekeeper's House, inserted
#    right after "Zum Schlafenden Auge / The Sleeping Eye" (row 12) and
#    before "Vielauge's Schloss 1 / Manyeyes' Castle 1" (old row 13).
# ---------------------------------------------------------------------
$maps = $wb.Worksheets.Item("Maps")
$maps.Rows("13:13").Insert()
$maps.Cells.Item(13, 1).Value = 376
$maps.Cells.Item(13, 2).Value = "Pf" + [char]0x00F6 + "rtnerhaus / Gatekeeper's House"
$maps.Cells.Item(13, 3).Value = "2D"
$maps.Cells.Item(13, 4).Value = "3 houses in cavetown"

$maps.Activate()
$maps.Range("A14").Select()

# ---------------------------------------------------------------------
# 3. Chests sheet: just a cursor/selection move (no data change).
# ---------------------------------------------------------------------
$chests = $wb.Worksheets.Item("Chests")
$chests.Activate()
$chests.Range("E24").Select()

# ---------------------------------------------------------------------
# 4. Places sheet: just a cursor/selection move (no data change).
# ---------------------------------------------------------------------
$places = $wb.Worksheets.Item("Places")
$places.Activate()
$places.Range("C8").Select()

# ---------------------------------------------------------------------
# 5. Back to the Todo sheet, which stays the active tab, with the
#    selection moved to A3 (the now-empty cell).
# ---------------------------------------------------------------------
$todo.Activate()
$todo.Range("A3").Select()
